{"js": "const body = context.document.body;\n\n// --- Edit 1: \"5 Expected outcomes\" section --------------------------------\n// The source text already reads \"...some of these cities look believable...\"\n// (the old run/bookmark split around \"citie\"/\"s\" is invisible at the text\n// level), so no textual change is required there.\n\n// --- Edit 2: add a new sentence about the user study to the activity-plan\n// narrative paragraph, right after \"...conclude the thesis. \" -------------\nconst conclude = body.search(\"conclude the thesis. \", { matchCase: true });\nconclude.load(\"text\");\nawait context.sync();\n\nif (conclude.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for 'conclude the thesis. ', found \" +\n      conclude.items.length\n  );\n}\n\nconclude.items[0].insertText(\n  \"The user study will be conducted during the testing and data gathering part of the project.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// --- Edit 3: rename the \"Testing\" activity-plan list item to \"Testing and\n// data gathering\" (the longer label only needs a single tab before the\n// date range instead of two) ------------------------------------------------\nconst testingLabel = body.search(\"Testing \", { matchCase: true });\ntestingLabel.load(\"text\");\nawait context.sync();\n\nif (testingLabel.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for 'Testing ', found \" +\n      testingLabel.items.length\n  );\n}\n\ntestingLabel.items[0].insertText(\"and data gathering \", Word.InsertLocation.after);\nawait context.sync();\n\n// Remove one of the two tab characters that now separate the label from the\n// date range (longer labels in this list use a single tab).\nconst doubleTab = body.search(\"gathering \\t\", { matchCase: true });\ndoubleTab.load(\"text\");\nawait context.sync();\n\nif (doubleTab.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for 'gathering ' + tab, found \" +\n      doubleTab.items.length\n  );\n}\n\ndoubleTab.items[0].insertText(\"gathering \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1: \"5 Expected outcomes\" section ---------------------------------\n# The source text already reads \"...some of these cities look believable...\"\n# (the legacy run/bookmark split around \"citie\"/\"s\" is invisible at the text\n# level), so no textual change is required there.\n\n# --- Edit 2: add a new sentence about the user study to the activity-plan\n# narrative paragraph, right after \"...conclude the thesis. \" --------------\n$concludeRange = $d.Content\n$concludeRange.Find.ClearFormatting()\n$concludeRange.Find.MatchCase = $true\n$concludeRange.Find.MatchWholeWord = $false\n$foundConclude = $concludeRange.Find.Execute(\"conclude the thesis. \")\nif (-not $foundConclude) {\n    throw \"Could not find 'conclude the thesis. ' in the document.\"\n}\n$concludeRange.Collapse(0)  # wdCollapseEnd\n$concludeRange.Text = \"The user study will be conducted during the testing and data gathering part of the project.\"\n\n# --- Edit 3: rename the \"Testing\" activity-plan list item to \"Testing and\n# data gathering\" (the longer label only needs a single tab before the\n# date range instead of two) -------------------------------------------------\n$testingRange = $d.Content\n$testingRange.Find.ClearFormatting()\n$testingRange.Find.MatchCase = $true\n$testingRange.Find.MatchWholeWord = $false\n$foundTesting = $testingRange.Find.Execute(\"Testing \")\nif (-not $foundTesting) {\n    throw \"Could not find 'Testing ' in the document.\"\n}\n$testingRange.Collapse(0)  # wdCollapseEnd\n$testingRange.Text = \"and data gathering \"\n\n# Remove one of the two tab characters that now separate the label from the\n# date range (longer labels in this list use a single tab).\n$tabRange = $d.Content\n$tabRange.Find.ClearFormatting()\n$tabRange.Find.MatchCase = $true\n$tabRange.Find.MatchWholeWord = $false\n$foundTab = $tabRange.Find.Execute(\"gathering \" + [char]9)\nif (-not $foundTab) {\n    throw \"Could not find 'gathering ' followed by a tab in the document.\"\n}\n$tabRange.Text = \"gathering \"\n"}
